$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the dates in C10 and C11
$ws.Range("C10").Value = 46003
$ws.Range("C11").Value = 45998

# Append a new row (row 38) with the new debtor record
$ws.Range("A38").Value = 37
$ws.Range("B38").Value = "ABC"
$ws.Range("C38").Value = 46006
$ws.Range("C38").NumberFormat = $ws.Range("C37").NumberFormat
$ws.Range("D38").Value = 789456
$ws.Range("E38").Value = $false
